# New weekly price observation for Berenjena @ Terminal Hortofrutícola Agro
# Chillán. The existing history (rows 68..103) is pushed down one row to make
# room, and the freshest record is written into the now-vacant row 68.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the whole historical block (A68:R103) down by one row, landing on
# A69:R104 — this carries the former last row (103) into the new row 104
# unchanged, and duplicates row 68's data into row 69 (about to be
# overwritten below with the real new values).
$historyBlock = $ws.Range("A68:R103").Value2
$ws.Range("A69:R104").Value2 = $historyBlock

# The bulk Value2 write above doesn't carry cell formatting, so the date
# column on the newly materialised row 104 needs its number format restored
# (matches every other row in column D).
$ws.Range("D104").NumberFormat = $ws.Range("D103").NumberFormat

# Write this week's new record into row 68 (Volumen/Calidad/Unidad/Origen/Kg
# o Unidades are unchanged from the prior entry at this slot; only the date
# and the three price columns plus the derived $/Kg move).
$ws.Range("D68").Value2 = 45134
$ws.Range("K68").Value2 = 8000
$ws.Range("L68").Value2 = 9000
$ws.Range("M68").Value2 = 8500
$ws.Range("P68").Value2 = 142
